# Apply BOM updates: add antenna connector and buck regulator rows
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mccoy_bom")

# Fill in new cell values in the same order the original author entered
# them, so that the shared-strings table comes out in the matching order.
$ws.Range("H9").Value = "https://www.digikey.com/en/products/detail/ablic-u-s-a-inc/S-85S1AB18-I6T1U/9489538"
$ws.Range("B9").Value = "IC REG BUCK 1.8V 200MA SNT-6A"
$ws.Range("B10").Value = "IC REG BUCK 1.1V 200MA SNT-6A"
$ws.Range("H10").Value = "https://www.digikey.com/en/products/detail/ablic-u-s-a-inc/S-85S1AB11-I6T1U/9489535"
$ws.Range("C9").Value = "U2"
$ws.Range("C10").Value = "U3"
$ws.Range("H11").Value = "https://www.digikey.com/en/products/detail/molex/0733910060/1465165"
$ws.Range("B11").Value = "CONN SMA RCPT STR 50 OHM PCB"

$ws.Range("I9").Value = 1
$ws.Range("I10").Value = 1
$ws.Range("I11").Value = 1

# Update selection to match author's final cursor position
$ws.Range("B3").Select()

# Update workbook window position
$excel.Windows.Item(1).Left = 18885
$excel.Windows.Item(1).Top = 7515
